# Scheduled data refresh: update leve crafting profit figures across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 20000.5
$ws.Cells.Item(17, 10).Value = 20000.5
$ws.Cells.Item(17, 12).Value = 60001.5
$ws.Cells.Item(17, 14).Value = -60337.5

$ws.Cells.Item(80, 8).Value = 208511.4
$ws.Cells.Item(80, 9).Value = 499
$ws.Cells.Item(80, 10).Value = 260514.5
$ws.Cells.Item(80, 11).Value = 1497
$ws.Cells.Item(80, 12).Value = 781543.5
$ws.Cells.Item(80, 13).Value = -499
$ws.Cells.Item(80, 14).Value = -783539.5

$ws.Cells.Item(83, 8).Value = 208511.4
$ws.Cells.Item(83, 9).Value = 499
$ws.Cells.Item(83, 10).Value = 260514.5
$ws.Cells.Item(83, 11).Value = 4491
$ws.Cells.Item(83, 12).Value = 2344630.5
$ws.Cells.Item(83, 13).Value = 501
$ws.Cells.Item(83, 14).Value = -2354614.5

$ws.Cells.Item(98, 8).Value = 3813.4666
$ws.Cells.Item(98, 9).Value = 899.5
$ws.Cells.Item(98, 11).Value = 899.5
$ws.Cells.Item(98, 13).Value = 598.5

$ws.Cells.Item(100, 8).Value = 1490.7142
$ws.Cells.Item(100, 9).Value = 787.4
$ws.Cells.Item(100, 10).Value = 3249
$ws.Cells.Item(100, 11).Value = 787.4
$ws.Cells.Item(100, 12).Value = 3249
$ws.Cells.Item(100, 13).Value = -246.4
$ws.Cells.Item(100, 14).Value = -4331

$ws.Cells.Item(116, 8).Value = 8713.5
$ws.Cells.Item(116, 9).Value = 9249.25
$ws.Cells.Item(116, 11).Value = 9249.25
$ws.Cells.Item(116, 13).Value = -5807.25

$ws.Cells.Item(122, 8).Value = 3813.4666
$ws.Cells.Item(122, 9).Value = 899.5
$ws.Cells.Item(122, 11).Value = 2698.5
$ws.Cells.Item(122, 13).Value = -248.5

$ws.Cells.Item(137, 8).Value = 2220.2
$ws.Cells.Item(137, 9).Value = 2220.2
$ws.Cells.Item(137, 11).Value = 6660.599999999999
$ws.Cells.Item(137, 13).Value = -4110.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(26, 8).Value = 1000
$ws.Cells.Item(26, 9).Value = 1000
$ws.Cells.Item(26, 11).Value = 1000
$ws.Cells.Item(26, 13).Value = -670

$ws.Cells.Item(122, 8).Value = 2158.5454
$ws.Cells.Item(122, 9).Value = 1765.8334
$ws.Cells.Item(122, 10).Value = 2629.8
$ws.Cells.Item(122, 11).Value = 5297.5002
$ws.Cells.Item(122, 12).Value = 7889.400000000001
$ws.Cells.Item(122, 13).Value = -2847.5002
$ws.Cells.Item(122, 14).Value = -12789.4

$ws.Cells.Item(132, 8).Value = 1747.4286
$ws.Cells.Item(132, 9).Value = 1705.5
$ws.Cells.Item(132, 10).Value = 1999
$ws.Cells.Item(132, 11).Value = 5116.5
$ws.Cells.Item(132, 12).Value = 5997
$ws.Cells.Item(132, 13).Value = -2586.5
$ws.Cells.Item(132, 14).Value = -11057

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1549.6666
$ws.Cells.Item(86, 9).Value = 1549.6666
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 1549.6666
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = ""
$ws.Cells.Item(86, 14).Value = -426.6666

$ws.Cells.Item(89, 8).Value = 1549.6666
$ws.Cells.Item(89, 9).Value = 1549.6666
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 7748.333000000001
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = ""
$ws.Cells.Item(89, 14).Value = -2132.333000000001

$ws.Cells.Item(99, 8).Value = 2180.7896
$ws.Cells.Item(99, 10).Value = 2471.3635
$ws.Cells.Item(99, 12).Value = 2471.3635
$ws.Cells.Item(99, 14).Value = -5467.363499999999

$ws.Cells.Item(134, 8).Value = 2629.1428
$ws.Cells.Item(134, 9).Value = 2900.8333
$ws.Cells.Item(134, 10).Value = 999
$ws.Cells.Item(134, 11).Value = 8702.499899999999
$ws.Cells.Item(134, 12).Value = 2997
$ws.Cells.Item(134, 13).Value = -6167.499899999999
$ws.Cells.Item(134, 14).Value = -8067

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2321.1177
$ws.Cells.Item(31, 9).Value = 1588.7778
$ws.Cells.Item(31, 11).Value = 1588.7778
$ws.Cells.Item(31, 13).Value = -1293.7778

$ws.Cells.Item(34, 8).Value = 2321.1177
$ws.Cells.Item(34, 9).Value = 1588.7778
$ws.Cells.Item(34, 11).Value = 1588.7778
$ws.Cells.Item(34, 13).Value = -1386.7778

$ws.Cells.Item(86, 8).Value = 20000
$ws.Cells.Item(86, 9).Value = 20000
$ws.Cells.Item(86, 11).Value = 20000
$ws.Cells.Item(86, 13).Value = -18877

$ws.Cells.Item(89, 8).Value = 20000
$ws.Cells.Item(89, 9).Value = 20000
$ws.Cells.Item(89, 11).Value = 100000
$ws.Cells.Item(89, 13).Value = -94384

$ws.Cells.Item(99, 8).Value = 2239.0833
$ws.Cells.Item(99, 9).Value = 1945.1666
$ws.Cells.Item(99, 11).Value = 1945.1666
$ws.Cells.Item(99, 13).Value = -447.1666

$ws.Cells.Item(105, 8).Value = 3024.9395
$ws.Cells.Item(105, 9).Value = 2220.1052
$ws.Cells.Item(105, 10).Value = 4117.2144
$ws.Cells.Item(105, 11).Value = 2220.1052
$ws.Cells.Item(105, 12).Value = 4117.2144
$ws.Cells.Item(105, 13).Value = -473.1052
$ws.Cells.Item(105, 14).Value = -7611.2144

$ws.Cells.Item(126, 8).Value = 2239.0833
$ws.Cells.Item(126, 9).Value = 1945.1666
$ws.Cells.Item(126, 11).Value = 5835.4998
$ws.Cells.Item(126, 13).Value = -3365.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(94, 8).Value = 10792.923
$ws.Cells.Item(94, 9).Value = 3115.4285
$ws.Cells.Item(94, 11).Value = 9346.2855
$ws.Cells.Item(94, 13).Value = -8670.2855

$ws.Cells.Item(98, 8).Value = 240.8
$ws.Cells.Item(98, 9).Value = 199
$ws.Cells.Item(98, 10).Value = 303.5
$ws.Cells.Item(98, 11).Value = 597
$ws.Cells.Item(98, 12).Value = 910.5
$ws.Cells.Item(98, 13).Value = 901
$ws.Cells.Item(98, 14).Value = -3906.5

$ws.Cells.Item(114, 8).Value = 3549.111
$ws.Cells.Item(114, 10).Value = 3158
$ws.Cells.Item(114, 12).Value = 9474
$ws.Cells.Item(114, 14).Value = -15982

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1672
$ws.Cells.Item(102, 9).Value = 1281.6666
$ws.Cells.Item(102, 10).Value = 3233.3333
$ws.Cells.Item(102, 11).Value = 1281.6666
$ws.Cells.Item(102, 12).Value = 3233.3333
$ws.Cells.Item(102, 13).Value = 340.3334
$ws.Cells.Item(102, 14).Value = -6477.3333

$ws.Cells.Item(122, 8).Value = 4138.375
$ws.Cells.Item(122, 9).Value = 2400.25
$ws.Cells.Item(122, 10).Value = 5876.5
$ws.Cells.Item(122, 11).Value = 7200.75
$ws.Cells.Item(122, 12).Value = 17629.5
$ws.Cells.Item(122, 13).Value = -4750.75
$ws.Cells.Item(122, 14).Value = -22529.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4381.4546
$ws.Cells.Item(40, 9).Value = 3750.5
$ws.Cells.Item(40, 11).Value = 3750.5
$ws.Cells.Item(40, 13).Value = -3614.5

$ws.Cells.Item(61, 8).Value = 2926.182
$ws.Cells.Item(61, 9).Value = 2351.3333
$ws.Cells.Item(61, 11).Value = 2351.3333
$ws.Cells.Item(61, 13).Value = -2149.3333

$ws.Cells.Item(113, 8).Value = 2926.182
$ws.Cells.Item(113, 9).Value = 2351.3333
$ws.Cells.Item(113, 11).Value = 2351.3333
$ws.Cells.Item(113, 13).Value = -181.3332999999998

$ws.Cells.Item(122, 8).Value = 6880.1333
$ws.Cells.Item(122, 9).Value = 7179.5264
$ws.Cells.Item(122, 11).Value = 21538.5792
$ws.Cells.Item(122, 13).Value = -19088.5792

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(98, 8).Value = 47999.5
$ws.Cells.Item(98, 10).Value = 47999.5
$ws.Cells.Item(98, 12).Value = 47999.5
$ws.Cells.Item(98, 14).Value = -53989.5

$ws.Cells.Item(123, 8).Value = 107933.336
$ws.Cells.Item(123, 10).Value = 107933.336
$ws.Cells.Item(123, 12).Value = 107933.336
$ws.Cells.Item(123, 14).Value = -117733.336

$ws.Cells.Item(136, 8).Value = 2258.077
$ws.Cells.Item(136, 9).Value = 2379.1667
$ws.Cells.Item(136, 10).Value = 805
$ws.Cells.Item(136, 11).Value = 7137.500100000001
$ws.Cells.Item(136, 12).Value = 2415
$ws.Cells.Item(136, 13).Value = -4587.500100000001
$ws.Cells.Item(136, 14).Value = -7515

$ws.Cells.Item(140, 8).Value = 79476
$ws.Cells.Item(140, 10).Value = 79476
$ws.Cells.Item(140, 12).Value = 79476
$ws.Cells.Item(140, 14).Value = -89836
